# Update the "想去人数" (want-to-go count) figures in column F
# for the "展览" and "全部类型" worksheets, which carry duplicate data.

$wb = $excel.ActiveWorkbook

$updates = @{
    6  = 45
    7  = 160
    9  = 31
    15 = 1059
    18 = 423
    23 = 1289
    24 = 2944
    27 = 757
    29 = 1643
    30 = 564
    32 = 29
    34 = 396
    38 = 16
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
